# Update the Scrum meeting attendance roll for the 8/29 meeting row (row 8).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark attendance "A" for Ryan Conyac (D), Michael McGregor (E), Yuchen Feng (F),
# Will McLain (G), Younouss Thiam (H).
$ws.Range("D8").Value = "A"
$ws.Range("E8").Value = "A"
$ws.Range("F8").Value = "A"
$ws.Range("G8").Value = "A"
$ws.Range("H8").Value = "A"

# Mark Brian Davis (I) as "T" (tardy).
$ws.Range("I8").Value = "T"

# Update the selection/view to reflect today's active cell.
$ws.Range("G8").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
